$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 753.9167
$ws.Range("J129").Value = 1398.75
$ws.Range("L129").Value = 4196.25
$ws.Range("N129").Value = -14196.25
$ws.Range("H137").Value = 13159303
$ws.Range("I137").Value = 1484.0667
$ws.Range("K137").Value = 4452.2001
$ws.Range("M137").Value = -1902.2001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1879.8
$ws.Range("I61").Value = 1428.5714
$ws.Range("J61").Value = 2932.6667
$ws.Range("K61").Value = 1428.5714
$ws.Range("L61").Value = 2932.6667
$ws.Range("M61").Value = -1216.5714
$ws.Range("N61").Value = -3356.6667
$ws.Range("H113").Value = 30398
$ws.Range("J113").Value = 30398
$ws.Range("L113").Value = 30398
$ws.Range("N113").Value = -39076
$ws.Range("H136").Value = 1879.8
$ws.Range("I136").Value = 1428.5714
$ws.Range("J136").Value = 2932.6667
$ws.Range("K136").Value = 4285.7142
$ws.Range("L136").Value = 8798.000100000001
$ws.Range("M136").Value = -1735.7142
$ws.Range("N136").Value = -13898.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1539.4559
$ws.Range("I31").Value = 1376.7931
$ws.Range("J31").Value = 1660.4103
$ws.Range("K31").Value = 1376.7931
$ws.Range("L31").Value = 1660.4103
$ws.Range("M31").Value = -1081.7931
$ws.Range("N31").Value = -2250.4103
$ws.Range("H34").Value = 1539.4559
$ws.Range("I34").Value = 1376.7931
$ws.Range("J34").Value = 1660.4103
$ws.Range("K34").Value = 1376.7931
$ws.Range("L34").Value = 1660.4103
$ws.Range("M34").Value = -1174.7931
$ws.Range("N34").Value = -2064.4103
$ws.Range("H86").Value = 2607.7942
$ws.Range("I86").Value = 2773.7917
$ws.Range("J86").Value = 2209.4
$ws.Range("K86").Value = 2773.7917
$ws.Range("L86").Value = 2209.4
$ws.Range("M86").Value = -1650.7917
$ws.Range("N86").Value = -4455.4
$ws.Range("H89").Value = 2607.7942
$ws.Range("I89").Value = 2773.7917
$ws.Range("J89").Value = 2209.4
$ws.Range("K89").Value = 13868.9585
$ws.Range("L89").Value = 11047
$ws.Range("M89").Value = -8252.958500000001
$ws.Range("N89").Value = -22279
$ws.Range("H99").Value = 2116.4285
$ws.Range("I99").Value = 2069.8667
$ws.Range("J99").Value = 2232.8333
$ws.Range("K99").Value = 2069.8667
$ws.Range("L99").Value = 2232.8333
$ws.Range("M99").Value = -571.8667
$ws.Range("N99").Value = -5228.8333
$ws.Range("H126").Value = 2116.4285
$ws.Range("I126").Value = 2069.8667
$ws.Range("J126").Value = 2232.8333
$ws.Range("K126").Value = 6209.6001
$ws.Range("L126").Value = 6698.499899999999
$ws.Range("M126").Value = -3739.6001
$ws.Range("N126").Value = -11638.4999
$ws.Range("H132").Value = 2156.72
$ws.Range("I132").Value = 1533.5883
$ws.Range("J132").Value = 3480.875
$ws.Range("K132").Value = 4600.7649
$ws.Range("L132").Value = 10442.625
$ws.Range("M132").Value = -2070.7649
$ws.Range("N132").Value = -15502.625
$ws.Range("H140").Value = 52700
$ws.Range("J140").Value = 52700
$ws.Range("L140").Value = 52700
$ws.Range("N140").Value = -63060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 57592530
$ws.Range("J9").Value = 40136932
$ws.Range("L9").Value = 120410796
$ws.Range("N9").Value = -120411244
$ws.Range("H22").Value = 21717172
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 21717172
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H49").Value = 933.3333
$ws.Range("I49").Value = 933.3333
$ws.Range("K49").Value = 2799.9999
$ws.Range("M49").Value = -2643.9999
$ws.Range("H117").Value = 2425.5833
$ws.Range("J117").Value = 3161.375
$ws.Range("L117").Value = 9484.125
$ws.Range("N117").Value = -16368.125
$ws.Range("H121").Value = 598.6
$ws.Range("I121").Value = 297.5
$ws.Range("J121").Value = 799.3333
$ws.Range("K121").Value = 892.5
$ws.Range("L121").Value = 2397.9999
$ws.Range("M121").Value = 417.5
$ws.Range("N121").Value = -5017.9999
$ws.Range("H126").Value = 4840
$ws.Range("I126").Value = 1933.3334
$ws.Range("J126").Value = 5566.6665
$ws.Range("K126").Value = 5800.0002
$ws.Range("L126").Value = 16699.9995
$ws.Range("M126").Value = -860.0002000000004
$ws.Range("N126").Value = -26579.9995
$ws.Range("H131").Value = 17997.08
$ws.Range("I131").Value = 336716.66
$ws.Range("J131").Value = 1791
$ws.Range("K131").Value = 1010149.98
$ws.Range("L131").Value = 5373
$ws.Range("M131").Value = -1005109.98
$ws.Range("N131").Value = -15453
$ws.Range("H139").Value = 791.53845
$ws.Range("I139").Value = 356.4
$ws.Range("J139").Value = 2242
$ws.Range("K139").Value = 1069.2
$ws.Range("L139").Value = 6726
$ws.Range("M139").Value = 4070.8
$ws.Range("N139").Value = -17006
$ws.Range("H141").Value = 10247.5
$ws.Range("I141").Value = 10015
$ws.Range("J141").Value = 10325
$ws.Range("K141").Value = 30045
$ws.Range("L141").Value = 30975
$ws.Range("M141").Value = -24865
$ws.Range("N141").Value = -41335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4500
$ws.Range("I80").Value = 4666.6665
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 4666.6665
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -3668.6665
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 4500
$ws.Range("I83").Value = 4666.6665
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 23333.3325
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -18341.3325
$ws.Range("N83").Value = -24984
$ws.Range("H132").Value = 2536.3823
$ws.Range("I132").Value = 1894.36
$ws.Range("K132").Value = 5683.08
$ws.Range("M132").Value = -3153.08

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1599.8
$ws.Range("I7").Value = 1666.3334
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 1666.3334
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -1554.3334
$ws.Range("N7").Value = -1724
$ws.Range("H16").Value = 503.5
$ws.Range("I16").Value = 499.27274
$ws.Range("J16").Value = 550
$ws.Range("K16").Value = 499.27274
$ws.Range("L16").Value = 550
$ws.Range("M16").Value = -329.27274
$ws.Range("N16").Value = -890
$ws.Range("H55").Value = 205.625
$ws.Range("I55").Value = 148.22223
$ws.Range("J55").Value = 279.42856
$ws.Range("K55").Value = 148.22223
$ws.Range("L55").Value = 279.42856
$ws.Range("M55").Value = 24.77777
$ws.Range("N55").Value = -625.4285600000001
$ws.Range("H109").Value = 14250
$ws.Range("J109").Value = 14250
$ws.Range("L109").Value = 14250
$ws.Range("N109").Value = -17024
$ws.Range("H122").Value = 2983.2812
$ws.Range("I122").Value = 2970.862
$ws.Range("K122").Value = 8912.585999999999
$ws.Range("M122").Value = -6462.585999999999
$ws.Range("H126").Value = 1599.8
$ws.Range("I126").Value = 1666.3334
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 4999.0002
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -2529.0002
$ws.Range("N126").Value = -9440
$ws.Range("H132").Value = 3790879.8
$ws.Range("J132").Value = 3856.2856
$ws.Range("L132").Value = 11568.8568
$ws.Range("N132").Value = -16628.8568
$ws.Range("H136").Value = 4133.75
$ws.Range("I136").Value = 1867.1428
$ws.Range("J136").Value = 20000
$ws.Range("K136").Value = 5601.428400000001
$ws.Range("L136").Value = 60000
$ws.Range("M136").Value = -3051.428400000001
$ws.Range("N136").Value = -65100
$ws.Range("H139").Value = 37350
$ws.Range("J139").Value = 37350
$ws.Range("L139").Value = 37350
$ws.Range("N139").Value = -47630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 950.7
$ws.Range("I81").Value = 900.7778
$ws.Range("J81").Value = 1400
$ws.Range("K81").Value = 1801.5556
$ws.Range("L81").Value = 2800
$ws.Range("M81").Value = -740.5555999999999
$ws.Range("N81").Value = -4922
$ws.Range("H84").Value = 950.7
$ws.Range("I84").Value = 900.7778
$ws.Range("J84").Value = 1400
$ws.Range("K84").Value = 9007.778
$ws.Range("L84").Value = 14000
$ws.Range("M84").Value = -3703.778
$ws.Range("N84").Value = -24608
$ws.Range("H140").Value = 48266.668
$ws.Range("J140").Value = 48266.668
$ws.Range("L140").Value = 48266.668
$ws.Range("N140").Value = -58626.668
